# Generate Report for Handback
# The a38ee76d-c2a0-4161-98f4-99e9098f354a.md file has just been handed back
# (target content is now in sync with en-US source). Update the status on the
# Overview sheet plus both language sheets, and stamp the "Latest Handback
# DateTime" for each language.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B3").Value = $statusText
$ws.Range("C3").Value = $statusText

# --- zh-cn sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("B3").Value = $statusText
$ws.Range("G3").Value = "2016-03-04 06:00:35"

# --- de-de sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("B3").Value = $statusText
$ws.Range("G3").Value = "2016-03-04 06:01:01"
